$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets.Item("autonomous_t26")
$ws.Delete()
$excel.DisplayAlerts = $true
